$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell "time_taken" in F1, matching the style of the
# other header cells (B1:E1) by copying formats from E1.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# time_taken values for each data row (F2:F11), plain/default style
$times = @(
    "2021-10-05 13:42:18.644424",
    "2021-10-05 13:42:18.644434",
    "2021-10-05 13:42:18.644437",
    "2021-10-05 13:42:18.644440",
    "2021-10-05 13:42:18.644443",
    "2021-10-05 13:42:18.644445",
    "2021-10-05 13:42:18.644448",
    "2021-10-05 13:42:18.644450",
    "2021-10-05 13:42:18.644453",
    "2021-10-05 13:42:18.644455"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
